# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Refreshed "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 09:22"

# Refreshed country stats. Several rows shifted by one position because a
# country's case count overtook its neighbours (Ucrania, Oman) or a
# duplicate-looking row got its real count (Anguila/Bonaire), so both the
# country name (column A) and the stats (columns B:H) need to be restated
# for every affected row.
$data = @(
    @(46, "Ucrania",               3102, 325, 97,  2912, 45, 10, 93),
    @(47, "Catar",                 2979, 0,   275, 2697, 37, 0,  7),
    @(48, "Finlandia",             2974, 0,   300, 2618, 80, 0,  56),
    @(49, "Republica Dominicana",  2967, 0,   131, 2663, 147,0,  173),
    @(60, "Moldavia",              1662, 0,   94,  1535, 80, 2,  33),
    @(80, "Oman",                  727,  128, 124, 599,  3,  0,  4),
    @(81, "Crucero",               712,  0,   619, 82,   10, 0,  11),
    @(82, "Tunez",                 707,  0,   43,  633,  85, 0,  31),
    @(83, "Bulgaria",              676,  1,   71,  574,  36, 2,  31),
    @(84, "Cuba",                  669,  0,   92,  559,  11, 0,  18),
    @(85, "Letonia",               653,  2,   16,  632,  2,  0,  5),
    @(86, "Principado de Andorra", 638,  0,   128, 481,  17, 0,  29),
    @(87, "Republica de Chipre",   633,  0,   65,  557,  8,  0,  11),
    @(88, "Libano",                630,  0,   80,  530,  34, 0,  20),
    @(89, "Banglades",             621,  0,   39,  548,  1,  0,  34),
    @(90, "Afganistan",            607,  0,   32,  556,  0,  1,  19),
    @(150, "Bahamas",              47,   1,   6,   33,   1,  0,  8),
    @(210, "Bonaire, San Eustaquio y Saba", 3, 0,  0,   3,    0,  0,  0),
    @(211, "Anguila",              3,    0,   1,   2,    0,  0,  0)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
}
